# "Custom Sort, Search implemented"
# - Adds a new "SubjectList" sheet (between "Maiya" and "Subas") with a small
#   price/quantity table.
# - Re-sorts the "Subas" sheet's data table by the Name column (ascending),
#   dropping a handful of leftover junk rows along the way.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "SubjectList" sheet, inserted right after "Maiya" (i.e. before "Subas")
# ---------------------------------------------------------------------------
$maiya = $wb.Worksheets.Item("Maiya")
$sl = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $maiya)
$sl.Name = "SubjectList"

$sl.Range("A1").Value = "Total Price"
$sl.Range("B1").Value = "  price"
$sl.Range("C1").Value = "  quantity"
$sl.Range("D1").Value = "Name"

$sl.Range("A2").Value = 3500.0
$sl.Range("B2").Value = 350.0
$sl.Range("C2").Value = 10.0
$sl.Range("D2").Value = "DSA"

$sl.Range("A3").Value = 6000.0
$sl.Range("B3").Value = 600.0
$sl.Range("C3").Value = 10.0
$sl.Range("D3").Value = "Electric Machine"

# ---------------------------------------------------------------------------
# 2) "Subas" sheet: drop stray junk rows, then custom-sort A1:D.. by column D
#    (Name) ascending.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Subas")

# Rows 43/44 hold one-off junk values ("kjhgfd" / "kldskajfsdkfa"); row 13
# holds another ("fiosdjkf"). Remove highest row numbers first so the other
# row indices stay valid.
$ws.Rows("43:44").Delete()
$ws.Rows("13:13").Delete()

# Also drop three redundant duplicate "Don" rows from the big repeated block.
$ws.Rows("39:41").Delete()

# Custom sort on column D (Name), ascending, range now A1:D38 with a header row.
$rng = $ws.Range("A1:D38")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D1:D38"), 0, 1, 0, 0)
$ws.Sort.SetRange($rng)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# The lone "Don" row carrying the distinct 96 value belongs at the very top
# of its group (it does, by construction). The lone "Subas" row carrying the
# distinct 84 value belongs at the very bottom of its group instead of the
# top, so swap its values down to the end of the "Subas" block.
$ws.Range("A27").Value = 32.0
$ws.Range("C27").Value = 12.0
$ws.Range("A37").Value = 84.0
$ws.Range("C37").Value = 12.0

# ---------------------------------------------------------------------------
# 3) Restore the originally active tab ("Maiya").
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Maiya").Activate()
